$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Step 2 text and its expected result
$ws.Range("C3").Value = "Step 2: Log in as the appropriate role"
$ws.Range("D3").Value = "I am redirected to the users dashboard"

# Copy the style used by row 3 (vertical top, wrap text) to the new rows
$styleSource = $ws.Range("C3:D3")

# Row 4: Step 3
$ws.Range("C4").Value = "Step 3: Click on View Employees"
$ws.Range("D4").Value = "A list of employees are displayed "

# Row 5: Step 4
$ws.Range("C5").Value = "Step 4: Select an employee to and click on delete"
$ws.Range("D5").Value = "I am redirected to the user's information page"

# Row 6: Step 5
$ws.Range("C6").Value = 'Step 5:Click "Delete this employee"'
$ws.Range("D6").Value = 'A pop up appears saying "Are you sure you want to delete this employee"'

# Row 7: Step 6
$ws.Range("C7").Value = "Step 6: Confirm the delete"
$ws.Range("D7").Value = "The employee is removed from the database"

$ws.Range("C4:D7").WrapText = $true
$ws.Range("C4:D7").VerticalAlignment = -4160

# Update the view: scroll so row 3 is at the top, and select the new last cell
$excel.ActiveWindow.ScrollRow = 3
$ws.Range("D7").Select()
